# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns,
# and swaps three pairs of coin rows (B/C/D/E) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.081.72"
$ws.Range("E2").Value = "'  -1.39%  "
$ws.Range("D3").Value = "'2.105.35"
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("E4").Value = "'  -0.69%  "
$ws.Range("D5").Value = "'349.99"
$ws.Range("E5").Value = "'  +4.31%  "
$ws.Range("E6").Value = "'  -0.59%  "
$ws.Range("D7").Value = "'0.5153"
$ws.Range("E7").Value = "'  -1.19%  "
$ws.Range("D8").Value = "'0.4466"
$ws.Range("E8").Value = "'  -1.36%  "
$ws.Range("D9").Value = "'52.64"
$ws.Range("E9").Value = "'  -4.39%  "
$ws.Range("D10").Value = "'0.08960"
$ws.Range("E10").Value = "'  -0.39%  "
$ws.Range("D11").Value = "'1.179"
$ws.Range("E11").Value = "'  +1.03%  "
$ws.Range("D12").Value = "'25.89"
$ws.Range("E12").Value = "'  +5.56%  "
$ws.Range("D13").Value = "'2.106.90"
$ws.Range("E13").Value = "'  -0.13%  "
$ws.Range("D14").Value = "'8.236"
$ws.Range("E14").Value = "'  +1.91%  "
$ws.Range("D15").Value = "'6.752"
$ws.Range("E15").Value = "'  -0.84%  "
$ws.Range("D16").Value = "'99.21"
$ws.Range("E16").Value = "'  +2.46%  "
$ws.Range("D17").Value = "'0.00001151"
$ws.Range("E17").Value = "'  -1.28%  "
$ws.Range("E18").Value = "'  -0.64%  "
$ws.Range("D19").Value = "'20.78"
$ws.Range("E19").Value = "'  +7.58%  "
$ws.Range("D20").Value = "'0.06669"
$ws.Range("E20").Value = "'  -0.10%  "
$ws.Range("E21").Value = "'  -0.64%  "
$ws.Range("D22").Value = "'6.262"
$ws.Range("E22").Value = "'  +0.72%  "
$ws.Range("D23").Value = "'30.168.35"
$ws.Range("E23").Value = "'  -1.32%  "
$ws.Range("D24").Value = "'12.89"
$ws.Range("E24").Value = "'  +0.46%  "
$ws.Range("D25").Value = "'2.352"
$ws.Range("E25").Value = "'  -0.22%  "
$ws.Range("D26").Value = "'2.357.05"
$ws.Range("E26").Value = "'  +0.03%  "
$ws.Range("E27").Value = "'  -0.80%  "
$ws.Range("E28").Value = "'  +2.35%  "
$ws.Range("D29").Value = "'162.67"
$ws.Range("E29").Value = "'  -0.34%  "
$ws.Range("D30").Value = "'133.71"
$ws.Range("E30").Value = "'  +0.34%  "
$ws.Range("D31").Value = "'1.182"
$ws.Range("E31").Value = "'  -2.25%  "
$ws.Range("E32").Value = "'  +0.43%  "
$ws.Range("D33").Value = "'1.640"
$ws.Range("E33").Value = "'  +0.68%  "
$ws.Range("D34").Value = "'6.272"
$ws.Range("E34").Value = "'  -0.80%  "
$ws.Range("D35").Value = "'3.968"
$ws.Range("E35").Value = "'  +0.36%  "
$ws.Range("B36").Value = "'FraxShare"
$ws.Range("C36").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'10.32"
$ws.Range("E36").Value = "'  -0.45%  "
$ws.Range("B37").Value = "'InternetComputer(DFINITY)"
$ws.Range("C37").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.934"
$ws.Range("E37").Value = "'  +1.13%  "
$ws.Range("D38").Value = "'0.02587"
$ws.Range("E38").Value = "'  -1.05%  "
$ws.Range("D39").Value = "'0.06859"
$ws.Range("E39").Value = "'  +1.05%  "
$ws.Range("B40").Value = "'Algorand"
$ws.Range("C40").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2314"
$ws.Range("E40").Value = "'  +0.43%  "
$ws.Range("B41").Value = "'Aptos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'12.76"
$ws.Range("E41").Value = "'  +1.85%  "
$ws.Range("E42").Value = "'  +0.40%  "
$ws.Range("D43").Value = "'1.257"
$ws.Range("E43").Value = "'  +0.47%  "
$ws.Range("D44").Value = "'14.39"
$ws.Range("E44").Value = "'  +1.99%  "
$ws.Range("D45").Value = "'2.340"
$ws.Range("E45").Value = "'  +1.82%  "
$ws.Range("D46").Value = "'0.6423"
$ws.Range("E46").Value = "'  +0.22%  "
$ws.Range("E47").Value = "'  +3.95%  "
$ws.Range("D48").Value = "'3.661"
$ws.Range("E48").Value = "'  -0.15%  "
$ws.Range("B49").Value = "'EOS"
$ws.Range("C49").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.225"
$ws.Range("E49").Value = "'  -1.77%  "
$ws.Range("B50").Value = "'Aave"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'83.45"
$ws.Range("E50").Value = "'  +0.69%  "
$ws.Range("D51").Value = "'0.07241"
$ws.Range("E51").Value = "'  +0.71%  "
